$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset rows got shuffled (dates/varieties/prices reassigned across
# rows 2-4 and 6-13; row 5 is untouched). Write the new values directly.

$rows = @{
    2  = @(44532, 'Brooks',  'Primera', 400, 27000, 28000, 27500, '$/bandeja 12 kilos', "Región de O'Higgins", 2292, 12)
    3  = @(44594, 'Santina', 'Primera', 160, 5000,  6000,  5500,  '$/bandeja 5 kilos',  "Región de O'Higgins", 1100, 5)
    4  = @(44210, 'Rainier', 'Segunda', 250, 21000, 22000, 21500, '$/caja 18 kilos',    "Región de O'Higgins", 1194, 18)
    6  = @(44229, 'Santina', 'Primera', 250, 6500,  7000,  6750,  '$/bandeja 5 kilos',  'Provincia de Curicó', 1350, 5)
    7  = @(44175, 'Rainier', 'Segunda', 270, 25000, 26000, 25500, '$/caja 18 kilos',    "Región de O'Higgins", 1417, 18)
    8  = @(44557, 'Lapins',  'Primera', 250, 9000,  10000, 9500,  '$/bandeja 10 kilos', 'Provincia de Curicó', 950,  10)
    9  = @(44571, 'Brooks',  'Segunda', 400, 8500,  9000,  8750,  '$/bandeja 10 kilos', "Región de O'Higgins", 875,  10)
    10 = @(44568, 'Santina', 'Segunda', 200, 15000, 16000, 15500, '$/bandeja 12 kilos', "Región de O'Higgins", 1292, 12)
    11 = @(44208, 'Lapins',  'Segunda', 200, 10500, 11000, 10750, '$/bandeja 12 kilos', 'Provincia de Curicó', 896,  12)
    12 = @(44537, 'Brooks',  'Primera', 200, 29000, 30000, 29500, '$/caja 20 kilos',    "Región de O'Higgins", 1475, 20)
    13 = @(44161, 'Bing',    'Primera', 160, 39000, 40000, 39500, '$/caja 20 kilos',    'Provincia de Curicó', 1975, 20)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D: Fecha
    $ws.Cells.Item($r, 11).Value = $vals[1]   # K: Variedad
    $ws.Cells.Item($r, 12).Value = $vals[2]   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $vals[3]   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals[4]   # N: Precio mínimo
    $ws.Cells.Item($r, 15).Value = $vals[5]   # O: Precio máximo
    $ws.Cells.Item($r, 16).Value = $vals[6]   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals[7]   # Q: Unidad de comercialización
    $ws.Cells.Item($r, 18).Value = $vals[8]   # R: Origen
    $ws.Cells.Item($r, 19).Value = $vals[9]   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals[10]  # T: Kg / unidad
}
